$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header field from "Folder" to "Functional Area"
$ws.Range("A1").Value = "Functional Area"

# Move the selection to A2, matching the final cursor position
$ws.Range("A2").Select()
